$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 689
    $ws.Range("F8").Value = 3375
    $ws.Range("F9").Value = 4279
    $ws.Range("F10").Value = 126
}
